# Insert 3 new rows of daily price data for "Plátano" (Vega Modelo de Temuco)
# right before the existing row 828, shifting the remaining rows down.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("828:830").Insert()

# Common (constant) values shared by every data row in this sheet
$mercadoId = 10
$mercado = "Vega Modelo de Temuco"
$region = "La Araucanía"
$codreg = 9
$tipo = "Fruta"
$productoId = 100108
$producto = "Tropicales y subtropicales"
$categoriaId = 100108006
$categoria = "Plátano"
$unidad = "$/caja 20 kilos"
$origen = "Ecuador"
$kgUnidad = 20

function Set-PlatanoRow($row, $fecha, $variedad, $calidad, $volumen, $precioMin, $precioMax, $precioProm, $precioKg) {
    $ws.Cells.Item($row, 1).Value = $mercadoId
    $ws.Cells.Item($row, 2).Value = $mercado
    $ws.Cells.Item($row, 3).Value = $region
    $ws.Cells.Item($row, 4).Value = $fecha
    $ws.Cells.Item($row, 5).Value = $codreg
    $ws.Cells.Item($row, 6).Value = $tipo
    $ws.Cells.Item($row, 7).Value = $productoId
    $ws.Cells.Item($row, 8).Value = $producto
    $ws.Cells.Item($row, 9).Value = $categoriaId
    $ws.Cells.Item($row, 10).Value = $categoria
    $ws.Cells.Item($row, 11).Value = $variedad
    $ws.Cells.Item($row, 12).Value = $calidad
    $ws.Cells.Item($row, 13).Value = $volumen
    $ws.Cells.Item($row, 14).Value = $precioMin
    $ws.Cells.Item($row, 15).Value = $precioMax
    $ws.Cells.Item($row, 16).Value = $precioProm
    $ws.Cells.Item($row, 17).Value = $unidad
    $ws.Cells.Item($row, 18).Value = $origen
    $ws.Cells.Item($row, 19).Value = $precioKg
    $ws.Cells.Item($row, 20).Value = $kgUnidad
}

# New row 828: Barraganete / Maduro
Set-PlatanoRow 828 44946 "Barraganete" "Maduro" 556 40000 40000 40000 2000

# New row 829: Barraganete / Verde
Set-PlatanoRow 829 44946 "Barraganete" "Verde" 50 35000 35000 35000 1750

# New row 830: Sin especificar / Pintón
Set-PlatanoRow 830 44946 "Sin especificar" "Pintón" 1500 24000 26000 24733 1237
